$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column, matching the style of the existing header row (copy
# formatting from G1, the last existing header cell, then set the new value)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for the data rows
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
